$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

$ws.Cells.Item(2, 2).Value = 'Bitcoin'
$ws.Cells.Item(2, 3).Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Cells.Item(2, 4).Value = '27.476.84'
$ws.Cells.Item(2, 5).Value = '  -5.41%  '
$ws.Cells.Item(3, 2).Value = 'Ethereum'
$ws.Cells.Item(3, 3).Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Cells.Item(3, 4).Value = '1.839.49'
$ws.Cells.Item(3, 5).Value = '  -4.40%  '
$ws.Cells.Item(4, 2).Value = 'TetherUSD'
$ws.Cells.Item(4, 3).Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextCell 4 4 '0.9995'
$ws.Cells.Item(4, 5).Value = '  -0.48%  '
$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell 5 4 '312.91'
$ws.Cells.Item(5, 5).Value = '  -3.91%  '
$ws.Cells.Item(6, 2).Value = 'USDC'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextCell 6 4 '0.9996'
$ws.Cells.Item(6, 5).Value = '  -0.35%  '
$ws.Cells.Item(7, 2).Value = 'XRP'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextCell 7 4 '0.4226'
$ws.Cells.Item(7, 5).Value = '  -7.88%  '
$ws.Cells.Item(8, 2).Value = 'Cardano'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell 8 4 '0.3628'
$ws.Cells.Item(8, 5).Value = '  -5.06%  '
$ws.Cells.Item(9, 2).Value = 'OKB'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 9 4 '44.07'
$ws.Cells.Item(9, 5).Value = '  -3.72%  '
$ws.Cells.Item(10, 2).Value = 'Dogecoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 10 4 '0.07229'
$ws.Cells.Item(10, 5).Value = '  -6.79%  '
$ws.Cells.Item(11, 2).Value = 'Polygon'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 11 4 '0.9014'
$ws.Cells.Item(11, 5).Value = '  -7.98%  '
$ws.Cells.Item(12, 2).Value = 'Solana'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 12 4 '20.54'
$ws.Cells.Item(12, 5).Value = '  -9.24%  '
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.821.89'
$ws.Cells.Item(13, 5).Value = '  -5.82%  '
$ws.Cells.Item(14, 2).Value = 'Chainlink'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 14 4 '6.573'
$ws.Cells.Item(14, 5).Value = '  -5.55%  '
$ws.Cells.Item(15, 2).Value = 'Polkadot'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 15 4 '5.313'
$ws.Cells.Item(15, 5).Value = '  -6.85%  '
$ws.Cells.Item(16, 2).Value = 'TRON'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 16 4 '0.06813'
$ws.Cells.Item(16, 5).Value = '  -2.84%  '
$ws.Cells.Item(17, 2).Value = 'BinanceUSD'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 17 4 '1.001'
$ws.Cells.Item(17, 5).Value = '  -0.51%  '
$ws.Cells.Item(18, 2).Value = 'Litecoin'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 18 4 '77.34'
$ws.Cells.Item(18, 5).Value = '  -8.85%  '
$ws.Cells.Item(19, 2).Value = 'ShibaInu'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 19 4 '0.000008966'
$ws.Cells.Item(19, 5).Value = '  -5.32%  '
$ws.Cells.Item(20, 2).Value = 'Dai'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 20 4 '1.001'
$ws.Cells.Item(20, 5).Value = '  -0.21%  '
$ws.Cells.Item(21, 2).Value = 'Avalanche'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 21 4 '15.35'
$ws.Cells.Item(21, 5).Value = '  -8.14%  '
$ws.Cells.Item(22, 2).Value = 'WrappedBTC'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(22, 4).Value = '27.473.23'
$ws.Cells.Item(22, 5).Value = '  -5.49%  '
$ws.Cells.Item(23, 2).Value = 'Uniswap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 23 4 '4.944'
$ws.Cells.Item(23, 5).Value = '  -7.62%  '
$ws.Cells.Item(24, 2).Value = 'Cosmos'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 24 4 '10.63'
$ws.Cells.Item(24, 5).Value = '  -3.70%  '
$ws.Cells.Item(25, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(25, 4).Value = '2.090.48'
$ws.Cells.Item(25, 5).Value = '  -3.27%  '
$ws.Cells.Item(26, 2).Value = 'Toncoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 26 4 '2.027'
$ws.Cells.Item(26, 5).Value = '  -1.76%  '
$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 27 4 '152.07'
$ws.Cells.Item(27, 5).Value = '  -3.70%  '
$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 28 4 '18.17'
$ws.Cells.Item(28, 5).Value = '  -4.38%  '
$ws.Cells.Item(29, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 29 4 '5.294'
$ws.Cells.Item(29, 5).Value = '  -5.55%  '
$ws.Cells.Item(30, 2).Value = 'BitcoinCash'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 30 4 '110.90'
$ws.Cells.Item(30, 5).Value = '  -5.70%  '
$ws.Cells.Item(31, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 31 4 '1.685'
$ws.Cells.Item(31, 5).Value = '  -8.09%  '
$ws.Cells.Item(32, 2).Value = 'Stellar'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 32 4 '0.08849'
$ws.Cells.Item(32, 5).Value = '  -4.85%  '
$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 33 4 '0.7764'
$ws.Cells.Item(33, 5).Value = '  -9.66%  '
$ws.Cells.Item(34, 2).Value = 'Filecoin'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 34 4 '4.507'
$ws.Cells.Item(34, 5).Value = '  -11.54%  '
$ws.Cells.Item(35, 2).Value = 'HuobiToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 35 4 '2.857'
$ws.Cells.Item(35, 5).Value = '  -5.21%  '
$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 36 4 '1.076'
$ws.Cells.Item(36, 5).Value = '  -13.43%  '
$ws.Cells.Item(37, 2).Value = 'Frax'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextCell 37 4 '0.9998'
$ws.Cells.Item(37, 5).Value = '  -0.35%  '
$ws.Cells.Item(38, 2).Value = 'Hedera'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 38 4 '0.05375'
$ws.Cells.Item(38, 5).Value = '  -5.55%  '
$ws.Cells.Item(39, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 39 4 '1.093'
$ws.Cells.Item(39, 5).Value = '  -5.02%  '
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 40 4 '0.01935'
$ws.Cells.Item(40, 5).Value = '  -5.37%  '
$ws.Cells.Item(41, 2).Value = 'MXToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 41 4 '2.935'
$ws.Cells.Item(41, 5).Value = '  -4.69%  '
$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 42 4 '6.831'
$ws.Cells.Item(42, 5).Value = '  -8.08%  '
$ws.Cells.Item(43, 2).Value = 'TheSandbox'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 43 4 '0.5052'
$ws.Cells.Item(43, 5).Value = '  -8.16%  '
$ws.Cells.Item(44, 2).Value = 'Algorand'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 44 4 '0.1634'
$ws.Cells.Item(44, 5).Value = '  -6.77%  '
$ws.Cells.Item(45, 2).Value = 'Cronos'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 45 4 '0.06621'
$ws.Cells.Item(45, 5).Value = '  -4.55%  '
$ws.Cells.Item(46, 2).Value = 'Aptos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 46 4 '8.220'
$ws.Cells.Item(46, 5).Value = '  -12.15%  '
$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 47 4 '0.4718'
$ws.Cells.Item(47, 5).Value = '  -8.81%  '
$ws.Cells.Item(48, 2).Value = 'Quant'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 48 4 '105.18'
$ws.Cells.Item(48, 5).Value = '  -4.82%  '
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 49 4 '10.20'
$ws.Cells.Item(49, 5).Value = '  -8.82%  '
$ws.Cells.Item(50, 2).Value = 'PaxDollar'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 50 4 '0.9986'
$ws.Cells.Item(50, 5).Value = '  -0.46%  '
$ws.Cells.Item(51, 2).Value = 'NEARProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 51 4 '1.640'
$ws.Cells.Item(51, 5).Value = '  -6.71%  '
